# RMA Details Maintenance Grid sheet: replace the three sample RMA
# records (RMA-AQZZ-*) with a newly generated batch (RMA-Q7ZG-*),
# matching the "Sales Order Line" (E), "Shipper Line" (F) and "Id" (J)
# columns for rows 2-4.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

$ws.Range("E2").Value = "RMA-Q7ZG-001"
$ws.Range("F2").Value = "RMA-Q7ZG-1-1"
$ws.Range("J2").Value = "a7s5f000000xKBlAAM"

$ws.Range("E3").Value = "RMA-Q7ZG-002"
$ws.Range("F3").Value = "RMA-Q7ZG-1-2"
$ws.Range("J3").Value = "a7s5f000000xKBmAAM"

$ws.Range("E4").Value = "RMA-Q7ZG-003"
$ws.Range("F4").Value = "RMA-Q7ZG-1-3"
$ws.Range("J4").Value = "a7s5f000000xKBnAAM"
